$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 16 de Julio de 2020 a las 10:41'
$ws.Cells.Item(4, 2).Value = 3617408
$ws.Cells.Item(4, 3).Value = 581
$ws.Cells.Item(4, 4).Value = 1646128
$ws.Cells.Item(4, 5).Value = 1831130
$ws.Cells.Item(6, 2).Value = 972144
$ws.Cells.Item(6, 3).Value = 1975
$ws.Cells.Item(6, 4).Value = 613881
$ws.Cells.Item(6, 5).Value = 333327
$ws.Cells.Item(6, 7).Value = 7
$ws.Cells.Item(6, 8).Value = 24936
$ws.Cells.Item(7, 2).Value = 752797
$ws.Cells.Item(7, 3).Value = 6428
$ws.Cells.Item(7, 4).Value = 531692
$ws.Cells.Item(7, 5).Value = 209168
$ws.Cells.Item(7, 7).Value = 167
$ws.Cells.Item(7, 8).Value = 11937
$ws.Cells.Item(20, 2).Value = 196323
$ws.Cells.Item(20, 3).Value = 2733
$ws.Cells.Item(20, 4).Value = 106963
$ws.Cells.Item(20, 5).Value = 86864
$ws.Cells.Item(20, 7).Value = 39
$ws.Cells.Item(20, 8).Value = 2496
$ws.Cells.Item(45, 2).Value = 47126
$ws.Cells.Item(45, 3).Value = 248
$ws.Cells.Item(45, 5).Value = 4111
$ws.Cells.Item(46, 2).Value = 44714
$ws.Cells.Item(46, 3).Value = 526
$ws.Cells.Item(46, 4).Value = 20100
$ws.Cells.Item(46, 5).Value = 24234
$ws.Cells.Item(46, 7).Value = 4
$ws.Cells.Item(46, 8).Value = 380
$ws.Cells.Item(47, 2).Value = 39054
$ws.Cells.Item(47, 3).Value = 333
$ws.Cells.Item(47, 4).Value = 28928
$ws.Cells.Item(47, 5).Value = 8521
$ws.Cells.Item(47, 7).Value = 11
$ws.Cells.Item(47, 8).Value = 1605
$ws.Cells.Item(63, 2).Value = 19270
$ws.Cells.Item(63, 3).Value = 116
$ws.Cells.Item(63, 4).Value = 17244
$ws.Cells.Item(63, 5).Value = 1315
$ws.Cells.Item(63, 7).Value = 1
$ws.Cells.Item(63, 8).Value = 711
$ws.Cells.Item(74, 1).Value = 'El Salvador'
$ws.Cells.Item(74, 2).Value = 10957
$ws.Cells.Item(74, 3).Value = 312
$ws.Cells.Item(74, 4).Value = 6257
$ws.Cells.Item(74, 5).Value = 4402
$ws.Cells.Item(74, 7).Value = 12
$ws.Cells.Item(74, 8).Value = 298
$ws.Cells.Item(75, 1).Value = 'Australia'
$ws.Cells.Item(75, 2).Value = 10810
$ws.Cells.Item(75, 3).Value = 323
$ws.Cells.Item(75, 4).Value = 8035
$ws.Cells.Item(75, 5).Value = 2662
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = 113
$ws.Cells.Item(88, 5).Value = 5934
$ws.Cells.Item(88, 7).Value = 2
$ws.Cells.Item(88, 8).Value = 46
$ws.Cells.Item(110, 4).Value = 2007
$ws.Cells.Item(110, 5).Value = 656
$ws.Cells.Item(117, 4).Value = 1904
$ws.Cells.Item(117, 5).Value = 43
$ws.Cells.Item(118, 2).Value = 1951
$ws.Cells.Item(118, 3).Value = 24
$ws.Cells.Item(118, 4).Value = 1514
$ws.Cells.Item(118, 5).Value = 409
$ws.Cells.Item(120, 1).Value = 'Lituania'
$ws.Cells.Item(120, 2).Value = 1902
$ws.Cells.Item(120, 3).Value = 6
$ws.Cells.Item(120, 4).Value = 1593
$ws.Cells.Item(120, 5).Value = 230
$ws.Cells.Item(120, 8).Value = 79
$ws.Cells.Item(121, 1).Value = 'Zambia'
$ws.Cells.Item(121, 2).Value = 1895
$ws.Cells.Item(121, 4).Value = 1412
$ws.Cells.Item(121, 5).Value = 441
$ws.Cells.Item(121, 8).Value = 42
$ws.Cells.Item(124, 2).Value = 1837
$ws.Cells.Item(124, 3).Value = 57
$ws.Cells.Item(124, 4).Value = 892
$ws.Cells.Item(124, 5).Value = 926
$ws.Cells.Item(160, 2).Value = 452
$ws.Cells.Item(160, 3).Value = 1
$ws.Cells.Item(160, 5).Value = 5
$ws.Cells.Item(209, 1).Value = 'Groenlandia'
$ws.Cells.Item(210, 1).Value = 'Islas Malvinas'
